$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '90.613.18'
$ws.Range("E2").Value = '  +2.99%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.151.66'
$ws.Range("E3").Value = '  +4.03%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.81'
$ws.Range("E5").Value = '  +3.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '626.77'
$ws.Range("E6").Value = '  +3.52%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.16'
$ws.Range("E7").Value = '  +32.71%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.366'
$ws.Range("E8").Value = '  +2.57%  '

$ws.Range("E9").Value = '  -0.09%  '

$ws.Range("B10").Value = 'LidoStakedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.148.48'
$ws.Range("E10").Value = '  +3.92%  '

$ws.Range("B11").Value = 'Cardano'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.797'
$ws.Range("E11").Value = '  +23.85%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.200'
$ws.Range("E12").Value = '  +7.85%  '

$ws.Range("E13").Value = '  +6.54%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000244'
$ws.Range("E14").Value = '  +4.66%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '35.27'
$ws.Range("E15").Value = '  +11.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '90.465.76'
$ws.Range("E16").Value = '  +2.54%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.736.81'
$ws.Range("E17").Value = '  +3.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.175.92'
$ws.Range("E18").Value = '  +3.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.67'
$ws.Range("E19").Value = '  +10.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.27'
$ws.Range("E20").Value = '  +7.83%  '

$ws.Range("B21").Value = 'PEPE'
$ws.Range("C21").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000212'
$ws.Range("E21").Value = '  +4.12%  '

$ws.Range("B22").Value = 'BitcoinCash'
$ws.Range("C22").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '460.27'
$ws.Range("E22").Value = '  +10.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.97'
$ws.Range("E23").Value = '  +12.09%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.22'
$ws.Range("E24").Value = '  +6.88%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.92'
$ws.Range("E25").Value = '  +10.94%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '91.99'
$ws.Range("E26").Value = '  +10.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.15'
$ws.Range("E27").Value = '  +6.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.323.64'
$ws.Range("E28").Value = '  +3.70%  '

$ws.Range("E29").Value = '  +0.01%  '

$ws.Range("B30").Value = 'Cronos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.161'
$ws.Range("E30").Value = '  +0.99%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.20'
$ws.Range("E31").Value = '  +13.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -8.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.68'
$ws.Range("E33").Value = '  +19.54%  '

$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '517.31'
$ws.Range("E34").Value = '  +4.21%  '

$ws.Range("B35").Value = 'Stellar'
$ws.Range("C35").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.184'
$ws.Range("E35").Value = '  +33.99%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.66'
$ws.Range("E36").Value = '  +3.79%  '

$ws.Range("E37").Value = '  +13.26%  '

$ws.Range("E38").Value = '  +7.57%  '

$ws.Range("E39").Value = '  +5.60%  '

$ws.Range("E40").Value = '  +5.41%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0882'
$ws.Range("E41").Value = '  +28.79%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.21'
$ws.Range("E42").Value = '  +0.20%  '

$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.420'
$ws.Range("E43").Value = '  +16.76%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.09%  '

$ws.Range("E45").Value = '  +0.02%  '

$ws.Range("E46").Value = '  +7.80%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '147.47'
$ws.Range("E47").Value = '  +0.64%  '

$ws.Range("E48").Value = '  +13.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.51'
$ws.Range("E49").Value = '  +2.72%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.32'
$ws.Range("E50").Value = '  +11.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.665'
$ws.Range("E51").Value = '  +15.22%  '
